# "made few changes to linear regression"
#
# 1) Expand the cached date-field text from a 2-digit year to a 4-digit
#    year ("2/24/17" -> "2/24/2017") everywhere the datetimeFigureOut
#    field is cached: the slide master and every slide layout.
# 2) Fix the "Coeffcient" -> "Coefficient" typo on slide 5.

$p = $ppt.ActivePresentation

function Update-DateField {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "2/24/17") {
                $tr.Text = "2/24/2017"
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DateField $master.Shapes

# Every slide layout's date placeholder.
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Update-DateField $layout.Shapes
}

# Fix the misspelled "Coeffcient" on slide 5 ("Covariance and Correlation
# Coeffcient" -> "... Coefficient").
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTextFrame -and ($shp.TextFrame.TextRange.Text -like "*Coeffcient*")) {
        $bodyRange = $shp.TextFrame.TextRange
        for ($pp = 1; $pp -le $bodyRange.Paragraphs().Count; $pp++) {
            $para = $bodyRange.Paragraphs($pp)
            for ($rr = 1; $rr -le 3; $rr++) {
                $run = $para.Runs($rr)
                if ($run.Text -eq "Coeffcient") {
                    $run.Text = "Coefficient"
                }
            }
        }
    }
}
